# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the 3e694ef3-... file rows on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 20:36:22"
$wsZhCn.Range("H3").Value = "2016-03-20 20:36:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 20:36:25"
$wsDeDe.Range("H3").Value = "2016-03-20 20:36:49"
